$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "Itt rövid leírás olvasható a kategóriáról ("
$suffix = ")"

for ($r = 2; $r -le 333; $r++) {
    $idval = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $prefix + $idval + $suffix
}

$ws.Cells.Item(334, 1).Value = "FB"
$ws.Cells.Item(335, 1).Value = "RE"

$ws.Cells.Item(334, 2).Value = $prefix + "FB" + $suffix
$ws.Cells.Item(335, 2).Value = $prefix + "RE" + $suffix

$ws.Application.Goto($ws.Range("A335"), $true)
